# Updates the cryptos list worksheet with refreshed price / volume(1h) data.
# Rows 33/34 and 48/49 also had their coin ranking order swapped.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.487.78'
$ws.Range('E2').Value = '  -1.64%  '
$ws.Range('D3').Value = '2.634.86'
$ws.Range('E3').Value = '  -0.71%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '582.55'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.44%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '157.19'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.85%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.646'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +2.75%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('E9').Value = '  -3.39%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '5.82'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.24%  '
$ws.Range('E11').Value = '  -1.43%  '
$ws.Range('E12').Value = '  -0.13%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '28.68'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.56%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000188'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -4.60%  '
$ws.Range('D15').Value = '3.109.24'
$ws.Range('E15').Value = '  -0.67%  '
$ws.Range('D16').Value = '64.270.51'
$ws.Range('E16').Value = '  -1.78%  '
$ws.Range('D17').Value = '2.648.10'
$ws.Range('E17').Value = '  +0.06%  '
$ws.Range('E18').Value = '  -3.14%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.69'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.80%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.46'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.35%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '347.15'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.61%  '
$ws.Range('E22').Value = '  +0.05%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '68.20'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.14%  '
$ws.Range('E24').Value = '  +7.11%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.0000112'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.60%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.45'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.77%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '593.01'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +9.64%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.60'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.73%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.02'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.83%  '
$ws.Range('E30').Value = '  -1.42%  '
$ws.Range('E31').Value = '  +0.23%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.09'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.02%  '
$ws.Range('B33').Value = 'ImmutableX'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.74'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.30%  '
$ws.Range('B34').Value = 'RenderToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.71'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +4.51%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.35'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.60%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.414'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.54%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '20.05'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.62%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.997'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.20%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.94'
$ws.Range('D39').Style = 'Normal'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '154.54'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.87%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.44'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +5.75%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '158.40'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.73%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '4.02'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.59%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '23.45'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +4.07%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0604'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.30%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.634'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.32%  '
$ws.Range('B48').Value = 'VeChain'
$ws.Range('C48').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0255'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.19%  '
$ws.Range('B49').Value = 'Stellar'
$ws.Range('C49').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.102'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.27%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '19.22'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.17%  '
$ws.Range('D51').Value = '0.0₆0235'
$ws.Range('E51').Value = '  -6.32%  '
